$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.587.66"
$ws.Range("E2").Value = "  +2.63%  "
$ws.Range("D3").Value = "3.212.62"
$ws.Range("E3").Value = "  +2.01%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.61%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.210.37"
$ws.Range("E8").Value = "  +1.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.547"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.95%  "
$ws.Range("E10").Value = "  +3.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.15"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.04%  "
$ws.Range("E12").Value = "  +2.98%  "
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.63%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000255"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.67%  "
$ws.Range("D15").Value = "3.743.96"
$ws.Range("E15").Value = "  +2.06%  "
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.37%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "65.183.28"
$ws.Range("E18").Value = "  +2.33%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.215.28"
$ws.Range("E19").Value = "  +2.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "484.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.89%  "
$ws.Range("E21").Value = "  +5.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.774"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +12.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "83.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.71%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.83%  "
$ws.Range("E29").Value = "  +4.38%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.75%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.78%  "
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("E33").Value = "  +10.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.59%  "
$ws.Range("D35").Value = "0.0₃0902"
$ws.Range("E35").Value = "  +3.96%  "
$ws.Range("E36").Value = "  +4.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.35%  "
$ws.Range("E38").Value = "  +4.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "484.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "52.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.47%  "
$ws.Range("E42").Value = "  +8.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.303"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0386"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.05%  "
$ws.Range("D45").Value = "2.952.99"
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("E46").Value = "  +3.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.14%  "
$ws.Range("E49").Value = "  +7.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.62%  "
